# Update the "Recorded By" (column G) values on the "Session Analysis Results"
# sheet: reorder the comma-separated list of recorder identities so that
# entries are alphabetically ordered with the literal "System" entry moved
# to the end of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"   = "backup@backdoor.com, system, System"
    "G3"   = "dnasr281@gmail.com, System"
    "G5"   = "backup@backdoor.com, System"
    "G6"   = "dnasr281@gmail.com, System"
    "G7"   = "admin@admin.com, System"
    "G8"   = "backup@backdoor.com, System"
    "G28"  = "backup@backdoor.com, system, System"
    "G29"  = "dnasr281@gmail.com, System"
    "G31"  = "backup@backdoor.com, System"
    "G32"  = "dnasr281@gmail.com, System"
    "G33"  = "admin@admin.com, System"
    "G34"  = "backup@backdoor.com, System"
    "G54"  = "backup@backdoor.com, system, System"
    "G55"  = "dnasr281@gmail.com, System"
    "G57"  = "backup@backdoor.com, System"
    "G58"  = "dnasr281@gmail.com, System"
    "G59"  = "admin@admin.com, System"
    "G60"  = "backup@backdoor.com, System"
    "G80"  = "backup@backdoor.com, System"
    "G81"  = "backup@backdoor.com, System"
    "G82"  = "backup@backdoor.com, System"
    "G87"  = "admin@admin.com, dnasr281@gmail.com"
    "G106" = "backup@backdoor.com, System"
    "G107" = "backup@backdoor.com, System"
    "G108" = "backup@backdoor.com, System"
    "G113" = "admin@admin.com, dnasr281@gmail.com"
    "G132" = "backup@backdoor.com, System"
    "G133" = "backup@backdoor.com, System"
    "G134" = "backup@backdoor.com, System"
    "G139" = "admin@admin.com, dnasr281@gmail.com"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
